$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Period highlight selector (H2): 18 -> 27
$ws.Range("H2").Value = 27

# Row 15 - Develop System
$ws.Range("E15").Value = 12
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 1

# Row 16 - Write Dissertation
$ws.Range("E16").Value = 12
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 0.9

# Row 17 - Testing System
$ws.Range("E17").Value = 24
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1

# Row 18 - Documenting
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 1

# View changes: zoom and selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 55
$ws.Range("AG6").Select()

# Column width changes (values chosen so the stored OOXML "width" attribute
# lands as close as possible to the target given this runtime's internal
# character->pixel rounding, which is based on a Maximum Digit Width of 7px)
$ws.Columns.Item(7).ColumnWidth = 19.857142857142858
$ws.Columns.Item(8).ColumnWidth = 2.7142857142857144
